# DescData.xlsx: "unify the conception of DataNode, DataTable, Entity."
# The sheet that used to describe a generic "Property1" table now represents
# a DataNode, so rename the sheet accordingly and leave the cursor where the
# author of the change last left it (D37) after touching the header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: Property1 -> DataNode
$ws.Name = "DataNode"

# The two wrapped-text header rows (column titles @ row1, row8) are drawn
# one point shorter than before.
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# Leave the selection on D37, matching the cursor position saved with the
# workbook.
[void]$ws.Range("D37").Select()
